$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 247, pushing the existing rows (and all rows
# below them) down by two. This matches the dimension growing from
# A1:T269 to A1:T271.
$ws.Rows.Item(247).Insert()
$ws.Rows.Item(247).Insert()

# New row 247: "Especial" quality entry for Región de O'Higgins.
$ws.Range("A247").Value = 7
$ws.Range("B247").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C247").Value = "Ñuble"
$ws.Range("D247").Value = 45077
$ws.Range("E247").Value = 16
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100101
$ws.Range("H247").Value = "Berries"
$ws.Range("I247").Value = 100101007
$ws.Range("J247").Value = "Kiwi"
$ws.Range("K247").Value = "Hayward"
$ws.Range("L247").Value = "Especial"
$ws.Range("M247").Value = 70
$ws.Range("N247").Value = 12000
$ws.Range("O247").Value = 13000
$ws.Range("P247").Value = 12714
$ws.Range("Q247").Value = "`$/bandeja 18 kilos"
$ws.Range("R247").Value = "Región de O'Higgins"
$ws.Range("S247").Value = 706
$ws.Range("T247").Value = 18

# New row 248: "Primera" quality entry for Región de O'Higgins.
$ws.Range("A248").Value = 7
$ws.Range("B248").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C248").Value = "Ñuble"
$ws.Range("D248").Value = 45077
$ws.Range("E248").Value = 16
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100101
$ws.Range("H248").Value = "Berries"
$ws.Range("I248").Value = 100101007
$ws.Range("J248").Value = "Kiwi"
$ws.Range("K248").Value = "Hayward"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 40
$ws.Range("N248").Value = 10000
$ws.Range("O248").Value = 10000
$ws.Range("P248").Value = 10000
$ws.Range("Q248").Value = "`$/bandeja 18 kilos"
$ws.Range("R248").Value = "Región de O'Higgins"
$ws.Range("S248").Value = 556
$ws.Range("T248").Value = 18
